# Parallel-execution update:
#  - Keep "LeaveEntitlement" as the single-row smoke-test sheet (header + TC001).
#  - Move the full original data set (TC001..TC004) into a new "Sheet1" so the
#    multi-scenario rows can be driven separately / in parallel.
#  - Shorten the EmployeeName values from "Fiona Grace"/"Hannah Flores" down to
#    just "Fiona"/"Hannah" on both sheets.

$wb = $excel.ActiveWorkbook
$leave = $wb.Worksheets.Item("LeaveEntitlement")

# 1) Duplicate the original sheet (with all 4 data rows) and place it right
#    after "LeaveEntitlement"; this becomes the new "Sheet1".
$leave.Copy($null, $leave)
$dataSheet = $wb.Worksheets.Item("LeaveEntitlement (2)")
$dataSheet.Name = "Sheet1"

# 2) On the new "Sheet1", trim the employee names.
$dataSheet.Range("F2").Value = "Fiona"
$dataSheet.Range("F3").Value = "Hannah"
$dataSheet.Range("F5").Value = "Hannah"

# 3) On "LeaveEntitlement", drop rows 3-5 (only header + TC001 remain) and
#    shorten the employee name on the row that's left.
$leave.Rows("3:5").Delete()
$leave.Range("F2").Value = "Fiona"

# 4) Restore view state: "LeaveEntitlement" stays the active/visible tab with
#    its selection moved to H6; "Sheet1" ends up with the whole sheet selected.
$dataSheet.Activate()
$dataSheet.Cells.Select()
$leave.Activate()
$leave.Range("H6").Select()
